# B6-PowerPoint.pptx edit script
#
# 1) Three tables (on the slides that used to be Google-Slides-exported
#    slides 14, 15 and 16) get their custom "Table_0" table style
#    ({E6AD10A1-80E6-4B6B-AD93-075FA7ADBABE}) replaced with the built-in
#    PowerPoint table style "Light Style 3"
#    ({74709A3B-D6CF-4B55-BADE-05F5C1B5ACD5}).
#
# 2) The deck's two themes are swapped: the theme actually driving the
#    slide master/slides (the "Integral"/Red Violet palette) becomes the
#    plain "Office Theme" palette, while the notes-only theme becomes the
#    "Integral"/Red Violet palette.

$p = $ppt.ActivePresentation

# --- 1) Retarget the three tables to the built-in "Light Style 3" style ---
$targetStyle = "{74709A3B-D6CF-4B55-BADE-05F5C1B5ACD5}"
foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($shapeIndex = 1; $shapeIndex -le $slide.Shapes.Count; $shapeIndex++) {
        $shape = $slide.Shapes.Item($shapeIndex)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyle)
        }
    }
}

# --- 2) Swap the active theme's colour scheme to the "Office" palette ---
# (The slide master's theme is the one that actually renders; its colour
# scheme currently holds the "Integral"/Red Violet palette and needs to
# become the standard Office palette.)
$officePalette = @(
    0,         # dk1       000000
    16777215,  # lt1       FFFFFF
    6968388,   # dk2       44546A
    15132391,  # lt2       E7E6E6
    13998939,  # accent1   5B9BD5
    3243501,   # accent2   ED7D31
    10855845,  # accent3   A5A5A5
    49407,     # accent4   FFC000
    12874308,  # accent5   4472C4
    4697456,   # accent6   70AD47
    12673797,  # hlink     0563C1
    7491477    # folHlink  954F72
)

$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officePalette[$i - 1]
}
